# "made few changes in keyword approach"
#
# Rewrites the BDD-style "General" sheet so column C holds descriptive
# rich-text guidance (plain lead-in + a non-bold quoted keyword + a bold
# call-to-action) instead of the old short keyword phrases, updates the
# sample first/last name data, widens column C, and moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-KeywordCell {
    param($Cell, $Prefix, $Middle, $Suffix)

    $full = $Prefix + $Middle + $Suffix
    $Cell.Value = $full

    # Base cell style -> fontId 1 (bold Calibri), no alignment overrides.
    $Cell.Style = "Normal"
    $Cell.Font.Bold = $true

    $midStart = $Prefix.Length + 1
    $midLen = $Middle.Length
    $suffixStart = $Prefix.Length + $Middle.Length + 1
    $suffixLen = $Suffix.Length

    # Quoted keyword in the middle: explicit (non-bold) run.
    $Cell.Characters($midStart, $midLen).Font.Bold = $false
    # Trailing call-to-action: explicit bold run.
    $Cell.Characters($suffixStart, $suffixLen).Font.Bold = $true
}

Set-KeywordCell `
    $ws.Cells.Item(2, 3) `
    'While filling the form, navigate to' `
    ' "Testzen Labs Form" ' `
    'to proceed with registration.'

Set-KeywordCell `
    $ws.Cells.Item(3, 3) `
    'Please ensure you correctly enter' `
    ' "First Name" ' `
    'before moving to the next field.'

Set-KeywordCell `
    $ws.Cells.Item(4, 3) `
    'You should carefully enter' `
    ' "Last Name" ' `
    'so that it matches your official documents.'

Set-KeywordCell `
    $ws.Cells.Item(5, 3) `
    'Before proceeding further, make sure to enter' `
    ' "Phone Number" ' `
    'to receive OTP verification.'

Set-KeywordCell `
    $ws.Cells.Item(6, 3) `
    'In the form, select' `
    ' "Country" ' `
    'from the dropdown list to specify your nationality.'

Set-KeywordCell `
    $ws.Cells.Item(7, 3) `
    'To complete your application, kindly upload' `
    ' "Resume" ' `
    'in the specified format.'

Set-KeywordCell `
    $ws.Cells.Item(8, 3) `
    'For gender identification, check the' `
    ' "Male" ' `
    'option if applicable.'

# Row 9 stays plain text, but still picks up the bold "keyword" style.
$c9 = $ws.Cells.Item(9, 3)
$c9.Value = 'To enhance security, generate a random number for the pin code before submission.'
$c9.Style = "Normal"
$c9.Font.Bold = $true

# Sample data updates.
$ws.Cells.Item(3, 6).Value = "Madan"
$ws.Cells.Item(4, 6).Value = "Reddy"

# Column C needs to be much wider to fit the new descriptive text.
$ws.Columns.Item(3).ColumnWidth = 84.3

# Row 9 no longer needs the extra height the old 2-line label required.
$ws.Rows.Item(9).AutoFit()

# Move the active selection.
$ws.Range("C8").Select()

Write-Output "done"
